# Auto-generated: update cryptos price/volume table (and shifted rows 46-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.388.23'
$ws.Range('D3').Value = '1.841.79'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = '''239.21'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '''0.6263'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''0.07429'
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('D9').Value = '''0.2893'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = '''24.93'
$ws.Range('E10').Value = '  +2.13%  '
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('D12').Value = '1.821.77'
$ws.Range('E12').Value = '  -1.20%  '
$ws.Range('D13').Value = '''4.975'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = '''0.6744'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').Value = '''0.00001027'
$ws.Range('E15').Value = '  -2.11%  '
$ws.Range('D16').Value = '''81.79'
$ws.Range('E16').Value = '  -0.28%  '
$ws.Range('D17').Value = '''6.215'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '29.410.61'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '''234.65'
$ws.Range('E19').Value = '  +2.50%  '
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').Value = '''7.293'
$ws.Range('D23').Value = '''1.001'
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').Value = '''158.49'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = '''8.488'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').Value = '''0.1344'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').Value = '''17.31'
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('D28').Value = '''0.07283'
$ws.Range('E28').Value = '  +13.55%  '
$ws.Range('D29').Value = '''1.463'
$ws.Range('E29').Value = '  +4.36%  '
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').Value = '''4.041'
$ws.Range('E31').Value = '  -1.31%  '
$ws.Range('D32').Value = '''4.028'
$ws.Range('E32').Value = '  -1.50%  '
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = '''0.6974'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '''2.572'
$ws.Range('D37').Value = '''0.01842'
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('D38').Value = '''6.921'
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('E39').Value = '  -0.79%  '
$ws.Range('D40').Value = '1.232.57'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').Value = '''0.9604'
$ws.Range('E41').Value = '  +4.97%  '
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('D43').Value = '1.997.73'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').Value = '''101.03'
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').Value = '''65.45'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''1.718'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '''6.949'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''8.858'
$ws.Range('E48').Value = '  -1.19%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.1134'
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = '''0.3897'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05658'
$ws.Range('E51').Value = '  -0.65%  '
